$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.1621188509149718
$ws.Range("B1").Value = 0.16204842881005987
$ws.Range("A2").Value = -0.13994738937956885
$ws.Range("B2").Value = 0.13962190436666777
$ws.Range("A3").Value = -0.097841229898975257
$ws.Range("B3").Value = 0.097487606958537754
$ws.Range("A4").Value = -0.089487607061844443
$ws.Range("B4").Value = 0.089159728343254585
$ws.Range("A5").Value = -0.086159728401002944
$ws.Range("B5").Value = 0.085054826246497406
$ws.Range("A6").Value = 0.014116777823094395
$ws.Range("B6").Value = -0.014159910203366266
$ws.Range("A7").Value = 0.024159910058167533
$ws.Range("B7").Value = -0.024169615862472238
$ws.Range("A8").Value = 0.034169615719868318
$ws.Range("B8").Value = -0.034229529377608259
$ws.Range("A9").Value = 0.036229529324243614
$ws.Range("B9").Value = -0.036303580414142722
$ws.Range("A10").Value = -0.013640906516583229
$ws.Range("B10").Value = 0.013641161932135049
$ws.Range("A11").Value = -0.010641161989101811
$ws.Range("B11").Value = 0.010640788098656095
$ws.Range("A12").Value = -0.0071407881612453039
$ws.Range("B12").Value = 0.0071355221119850931
$ws.Range("A13").Value = -0.003635522176233863
$ws.Range("B13").Value = 0.0036349277354759124
$ws.Range("A14").Value = 0.0043650721534316261
$ws.Range("B14").Value = -0.0043708690414856477
$ws.Range("A15").Value = 0.0053708690032578943
$ws.Range("B15").Value = -0.0053843114815688153
$ws.Range("A16").Value = 0.0073843114335963023
$ws.Range("B16").Value = -0.0074615955802648237
$ws.Range("A17").Value = -0.0040028947180967833
$ws.Range("B17").Value = 0.0039999999328257374
$ws.Range("A18").Value = -0.016101039493921832
$ws.Range("B18").Value = 0.016090539803762738
$ws.Range("A19").Value = -0.012090539846648429
$ws.Range("B19").Value = 0.012015883184206988
$ws.Range("A20").Value = -0.0080158832305201599
$ws.Range("B20").Value = 0.0080055519057573576
$ws.Range("A21").Value = -0.0040055519525772354
$ws.Range("B21").Value = 0.0039999999527751129
$ws.Range("A22").Value = -0.037780674522966962
$ws.Range("B22").Value = 0.037563503716903313
$ws.Range("A23").Value = -0.011764393789523098
$ws.Range("B23").Value = 0.011593161591165746
$ws.Range("A24").Value = -0.020096546120915271
$ws.Range("B24").Value = 0.019999999770527133
$ws.Range("A25").Value = -0.09717160413232051
$ws.Range("B25").Value = 0.097052927068888906
$ws.Range("A26").Value = -0.094552927139629261
$ws.Range("B26").Value = 0.094398853211298928
$ws.Range("A27").Value = -0.0918988532865046
$ws.Range("B27").Value = 0.090978483190872428
$ws.Range("A28").Value = -0.0889784832836229
$ws.Range("B28").Value = 0.088348594054735941
$ws.Range("A29").Value = -0.081348594215503667
$ws.Range("B29").Value = 0.081166550185327146
$ws.Range("A30").Value = -0.021166550903167991
$ws.Range("B30").Value = 0.021021455585731541
$ws.Range("A31").Value = -0.014021455759090529
$ws.Range("B31").Value = 0.014000782822897406
$ws.Range("A32").Value = -0.0040007830280366363
$ws.Range("B32").Value = 0.0039999998568518436
